$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("URI schema")

# 1) Insert 5 new rows before row 17 (shifts existing rows 17.. down to 22..)
$ws.Rows("17:21").Insert()

# 2) Populate the 5 new "Concepten" relation rows (17-21)
$newRows = @(
    @{ Row = 17; B = "focus";       D = "http://data.test.pdok.nl/catalogus/dso/id/concept/{focus}" },
    @{ Row = 18; B = "hetzelfde";   D = "http://data.test.pdok.nl/catalogus/dso/id/concept/{hetzelfde}" },
    @{ Row = 19; B = "gerelateerd"; D = "http://data.test.pdok.nl/catalogus/dso/id/concept/{gerelateerd}" },
    @{ Row = 20; B = "brederdan";   D = "http://data.test.pdok.nl/catalogus/dso/id/concept/{brederdan}" },
    @{ Row = 21; B = "engerdan";    D = "http://data.test.pdok.nl/catalogus/dso/id/concept/{engerdan}" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = "Concepten"
    $ws.Cells.Item($r, 2).Value = $nr.B
    $ws.Cells.Item($r, 4).Value = $nr.D
}

# 3) Rebuild every hyperlink on the sheet so refs line up with the new row numbers.
#    (Deleting the Hyperlinks collection in bulk is the only reliable way to clear
#    them in this runtime - per-item Delete() is a no-op.)
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D5"), "http://purl.org/dc/terms/{eigenschap}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "http://purl.org/iso25964/skos-thes", "{eigenschap}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D12"), "http://localhost:8080/catalogus/dsoprogramma/id/begrip/{specialisatie}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D13"), "http://localhost:8080/catalogus/dsoprogramma/id/begrip/{generalisatie}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D14"), "http://localhost:8080/catalogus/dsoprogramma/id/begrip/{onderdeel}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D15"), "http://localhost:8080/catalogus/dsoprogramma/id/begrip/{bestaatuit}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D16"), "http://localhost:8080/catalogus/dsoprogramma/id/begrip/{betrekkingop}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D22"), "http://localhost:8080/catalogus/dso/id/collection/{collectie}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D23"), "http://localhost:8080/dsoprogramma/id/begrip/{begrip}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D9"), "http://localhost:8080/dsoprogramma/id/begrip/{begrip}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "http://www.w3.org/1999/02/22-rdf-syntax-ns", "{eigenschap}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "http://www.w3.org/2004/02/skos/core", "{klasse}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D11"), "http://localhost:8080/catalogus/dso/id/concept/{bron}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D28"), "http://purl.org/dc/dcmitype/{subklasse}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D27"), "http://localhost:8080/catalogus/dso/id/concept/{bron}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D25"), "http://localhost:8080/catalogus/dso/id/collection/{waardelijst}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D26"), "http://localhost:8080/catalogus/dso/id/collection/{collectie}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D30"), "http://data.test.pdok.nl/catalogus/dso/id/concept/{toeleidingsbegrip}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D31"), "http://localhost:8080/dsoprogramma/id/begrip/{begrip}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "http://xmlns.com/foaf/0.1/{eigenschap}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D10"), "http://data.test.pdok.nl/catalogus/dso/id/concept/{domein}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D24"), "http://data.test.pdok.nl/catalogus/dso/id/concept/{domein}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D17"), "http://data.test.pdok.nl/catalogus/dso/id/concept/{focus}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D18"), "http://data.test.pdok.nl/catalogus/dso/id/concept/{hetzelfde}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D19"), "http://data.test.pdok.nl/catalogus/dso/id/concept/{gerelateerd}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D20"), "http://data.test.pdok.nl/catalogus/dso/id/concept/{brederdan}") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D21"), "http://data.test.pdok.nl/catalogus/dso/id/concept/{engerdan}") | Out-Null

# 4) Selection matches the authored state (cursor parked at E11)
$ws.Range("E11").Select()
